$wb = $excel.ActiveWorkbook

# --- Controls sheet: update n_sims (B2) and n_sex (B5) ---
$controls = $wb.Worksheets.Item("Controls")
$controls.Range("B2").Value2 = 100
$controls.Range("B5").Value2 = 2

# --- Maturity_At_Age: add row 3 (sex 2), duplicate of row 2 but with C3 = 2 ---
$maturity = $wb.Worksheets.Item("Maturity_At_Age")
$maturity.Range("A3").Value2 = 1
$maturity.Range("B3").Value2 = $maturity.Range("B2").Value2
$maturity.Range("C3").Value2 = 2
$maturity.Range("D3:AG3").Value2 = $maturity.Range("D2:AG2").Value2
[void]$maturity.Range("C3").Select()

# --- Weight_At_Age: add row 3 (sex 2), duplicate of row 2 but with C3 = 2 ---
$weight = $wb.Worksheets.Item("Weight_At_Age")
$weight.Range("A3").Value2 = 1
$weight.Range("B3").Value2 = $weight.Range("B2").Value2
$weight.Range("C3").Value2 = 2
$weight.Range("D3:AG3").Value2 = $weight.Range("D2:AG2").Value2
[void]$weight.Range("C3").Select()

# --- Make Controls the active/selected sheet (was Recruitment_Mortality) ---
$controls.Activate()
